$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Zetfast Loop" summary block that lives in columns H:L (rows 4-9)
# shifts up by one row: the old header row (H5:L5, style 15) moves to
# H4:L4, the old label row (H6:L6, style 11) moves to H5:L5, and the old
# data rows (H7:L9, styles 6/10) move to H6:L8. The old last data row
# (H9:L9) is removed entirely, shrinking the sheet's used range from
# A1:L9 to A1:L8.

# --- Step 1: copy cell formatting (styles only) from the old position to
# the new one, top-to-bottom so every source cell's style is read before
# it gets overwritten itself.
$ws.Range("H5:L5").Copy() | Out-Null
$ws.Range("H4:L4").PasteSpecial(-4122) | Out-Null

$ws.Range("H6:L6").Copy() | Out-Null
$ws.Range("H5:L5").PasteSpecial(-4122) | Out-Null

$ws.Range("H7:J7").Copy() | Out-Null
$ws.Range("H6:J6").PasteSpecial(-4122) | Out-Null
$ws.Range("K7:L7").Copy() | Out-Null
$ws.Range("K6:L6").PasteSpecial(-4122) | Out-Null

$ws.Range("H8:J8").Copy() | Out-Null
$ws.Range("H7:J7").PasteSpecial(-4122) | Out-Null
$ws.Range("K8:L8").Copy() | Out-Null
$ws.Range("K7:L7").PasteSpecial(-4122) | Out-Null

$ws.Range("H9:J9").Copy() | Out-Null
$ws.Range("H8:J8").PasteSpecial(-4122) | Out-Null
$ws.Range("K9:L9").Copy() | Out-Null
$ws.Range("K8:L8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Step 2: write the shifted text content. Plain .Value writes never
# disturb a cell's style, so these are safe in any order.
$ws.Range("H4").Value = "Zetfast Loop"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""

$ws.Range("H5").Value = "Zetfast loop with devices"
$ws.Range("I5").Value = "Gallery Type"
$ws.Range("J5").Value = "Label"
$ws.Range("K5").Value = "Battery Standby"
$ws.Range("L5").Value = "Alarm Load"

$ws.Range("H6").Value = "XLM800-Zetfas"
$ws.Range("I6").Value = "Loops"
$ws.Range("J6").Value = "XLM800-Zetfas-C"

$ws.Range("H7").Value = "Heat Sensor W3-S1"
$ws.Range("I7").Value = "Detectors"
$ws.Range("J7").Value = "Heat Sensor ... - 1"

$ws.Range("H8").Value = "Ion Sensor IR3-S1"
$ws.Range("I8").Value = "Detectors"
$ws.Range("J8").Value = "Ion Sensor I... - 2"

# --- Step 3: write the shifted numeric content for K:L. Writing a plain
# number into a cell that was formatted with a quoted-text style (the
# style backing K7/L7/K8/L8/K9/L9) drops that style's quote-prefix flag,
# so re-assert the style from the still-untouched K9:L9 donor right after
# each numeric write (K9:L9 itself is cleared only in step 4, once it is
# no longer needed as a format donor).
$ws.Range("K6").Value = 0.09
$ws.Range("L6").Value = 0.09
$ws.Range("K9:L9").Copy() | Out-Null
$ws.Range("K6:L6").PasteSpecial(-4122) | Out-Null

$ws.Range("K7").Value = 0.2
$ws.Range("L7").Value = 0.2
$ws.Range("K9:L9").Copy() | Out-Null
$ws.Range("K7:L7").PasteSpecial(-4122) | Out-Null

$ws.Range("K8").Value = 0.2
$ws.Range("L8").Value = 0.2
$ws.Range("K9:L9").Copy() | Out-Null
$ws.Range("K8:L8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Step 4: drop the now-duplicated last row (old row 9), which used to
# only hold H9:L9, so the sheet's used range shrinks back to row 8.
$ws.Range("H9:L9").Clear() | Out-Null

# --- Step 5: restore the selection state shown in the saved workbook.
$ws.Range("H4:L8").Select() | Out-Null
